$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update corrected marks/total in the marksheet "Total" row
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 15
$ws.Range("E12").Value = "15/140"
